$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Update existing Abilities sheet (sheet1) ---

# Row 7 (Swap): castTime (B7) drops from 0.3 to 0
$ws1.Range("B7").Value = 0

# New ability rows: Copy, Taunt, Buff
# Columns: A=AbilityType B=castTime C=cooldown D=power E=effectDuration F=ticks G=objectSpeed
$ws1.Range("A8").Value = "Copy"
$ws1.Range("B8").Value = 0
$ws1.Range("C8").Value = 0
$ws1.Range("D8").Value = 0
$ws1.Range("E8").Value = 10
$ws1.Range("F8").Value = 0
$ws1.Range("G8").Value = 0

$ws1.Range("A9").Value = "Taunt"
$ws1.Range("B9").Value = 0
$ws1.Range("C9").Value = 10
$ws1.Range("D9").Value = 0
$ws1.Range("E9").Value = 3
$ws1.Range("F9").Value = 0
$ws1.Range("G9").Value = 0

$ws1.Range("A10").Value = "Buff"
$ws1.Range("B10").Value = 0
$ws1.Range("C10").Value = 5
$ws1.Range("D10").Value = 2
$ws1.Range("E10").Value = 5
$ws1.Range("F10").Value = 0
$ws1.Range("G10").Value = 0

# --- Add new Characters sheet (after Abilities) ---
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Characters"

$ws2.Range("A1").Value = "CharacterType"
$ws2.Range("B1").Value = "Health"
$ws2.Range("C1").Value = "StartStun"

$ws2.Range("A2").Value = "Boss"
$ws2.Range("B2").Formula = "=SUM(B3:B11)"
$ws2.Range("C2").Value = 5

$ws2.Range("A3").Value = "Poker"
$ws2.Range("B3").Value = 100
$ws2.Range("C3").Value = 0

$ws2.Range("A4").Value = "Stuner"
$ws2.Range("B4").Value = 100
$ws2.Range("C4").Value = 0

$ws2.Range("A5").Value = "Nuker"
$ws2.Range("B5").Value = 50
$ws2.Range("C5").Value = 0

$ws2.Range("A6").Value = "Healer"
$ws2.Range("B6").Value = 100
$ws2.Range("C6").Value = 0

$ws2.Range("A7").Value = "Swapper"
$ws2.Range("B7").Value = 100
$ws2.Range("C7").Value = 0

$ws2.Range("A8").Value = "Barrier"
$ws2.Range("B8").Value = 100
$ws2.Range("C8").Value = 0

$ws2.Range("A9").Value = "Ditto"
$ws2.Range("B9").Value = 100
$ws2.Range("C9").Value = 0

$ws2.Range("A10").Value = "Taunter"
$ws2.Range("B10").Value = 100
$ws2.Range("C10").Value = 0

$ws2.Range("A11").Value = "Buffer"
$ws2.Range("B11").Value = 100
$ws2.Range("C11").Value = 0

# Column widths to roughly mirror the Abilities sheet layout
$ws2.Columns.Item(1).ColumnWidth = 23
$ws2.Range("C1:H1").ColumnWidth = 8.6

# Page setup (metric/narrow margins, portrait) to mirror the target layout
$ws2.PageSetup.Orientation = 1
$ws2.PageSetup.LeftMargin = 36.8503937007874
$ws2.PageSetup.RightMargin = 36.8503937007874
$ws2.PageSetup.TopMargin = 56.69291338582677
$ws2.PageSetup.BottomMargin = 56.69291338582677
$ws2.PageSetup.HeaderMargin = 22.677165354330707
$ws2.PageSetup.FooterMargin = 22.677165354330707

# Create the Characters table
$lo = $ws2.ListObjects.Add(1, $ws2.Range("A1:C11"), $null, 1)
$lo.Name = "Tabela1"
$lo.TableStyle = "TableStyleMedium8"

# Selections to mirror the target view state
$ws2.Range("C11").Select() | Out-Null
$ws1.Activate() | Out-Null
$ws1.Range("E7").Select() | Out-Null
